$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: collapse the "Login Page / Registration Page / Landing Page"
# sub-bullets (paragraphs 4-12) down into a single paragraph that keeps
# the "Login Page" paragraph's list formatting (ListParagraph, ilvl 0,
# numId 4) but only the last item's run ("User database using hibernate").
# The _GoBack bookmark that used to sit on the final (empty) paragraph of
# the document moves here too, right before that surviving run.
# ------------------------------------------------------------------

$d.Bookmarks.Item('_GoBack').Delete()

$pLoginPage = $d.Paragraphs.Item(4)
$pHibernate = $d.Paragraphs.Item(13)
$killRange = $d.Range($pLoginPage.Range.Start, $pHibernate.Range.Start)
$killRange.Delete()

$pMerged = $d.Paragraphs.Item(4)
$bmTarget = $d.Range($pMerged.Range.Start, $pMerged.Range.Start)
$d.Bookmarks.Add('_GoBack', $bmTarget)

# ------------------------------------------------------------------
# Part 2: the last paragraph of the document (previously empty, only
# holding the _GoBack bookmark) gains a write-up about the servlet/jsp
# work, followed by two more new bullet paragraphs with the same list
# formatting (ListParagraph, numId 5, ilvl 0).
#
# To make sure the freshly typed runs pick up the surrounding run
# formatting (rFonts/sz/szCs) instead of coming out as plain/unformatted
# runs, the three new paragraphs are split off of the end of the
# preceding ("After reading through the requirements...") paragraph,
# which already carries a properly formatted run. That leaves one extra
# blank paragraph (the original bookmark-only one) trailing at the very
# end of the document, which is then merged away by deleting its
# preceding paragraph mark.
# ------------------------------------------------------------------

$pPrev = $d.Paragraphs.Item(11)
$splitPoint = $d.Range($pPrev.Range.End - 1, $pPrev.Range.End - 1)
$splitPoint.InsertAfter("`rIn order to link all the servlets and jsp files properly, I had to use a few different methods for each for the different syntaxes. After they were linked, I could start writing the forms for the log-in a register pages. After those I could start writing the servlet code in order to make the proper checks. I created a new database so that the servlets had a server to connect to, but as of right now they don" + [char]0x2019 + "t have any members.`rAfter I wrote some of the basic functions, I decided to make another jsp file to notify the user whenever they have an unsuccessful login. This page worked out very well, but I had some issues getting the login screen to have successful logins`rAfter I finished the issue with logging in, I worked on the register pages. Making sure the servlet was inserting data into the database and that there were no duplicate usernames. If someone does try to reuse a username, they are redirected to another page")

$pBeforeBlank = $d.Paragraphs.Item(14)
$trailingMark = $d.Range($pBeforeBlank.Range.End - 1, $pBeforeBlank.Range.End)
$trailingMark.Delete()

Write-Output "done: paragraphs=$($d.Paragraphs.Count)"
